$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Datatype DateWrapper" block (rows 14-18)
$ws.Range("B14").Value = "Datatype DateWrapper"

$ws.Range("B15").Value = "Date"
$ws.Range("C15").Value = "date1"

$ws.Range("B16").Value = "Date"
$ws.Range("C16").Value = "date2"

$ws.Range("B17").Value = "Date"
$ws.Range("C17").Value = "date3"

$ws.Range("B18").Value = "Date"
$ws.Range("C18").Value = "date4"

# New "Spreadsheet DateWrapper spr(DateWrapper dw)" test block (rows 23-25)
$ws.Range("B24").Value = "Steps"
$ws.Range("C24").Value = "Formula"

$ws.Range("B25").Value = "RETURN"
$ws.Range("C25").Value = "'= dw"

$ws.Range("B23").Value = "Spreadsheet DateWrapper spr(DateWrapper dw)"

# Center the merged header above the DateWrapper fields and merge it
$ws.Range("B14:C14").HorizontalAlignment = -4108
$ws.Range("B14:C14").Merge($false)

# Restore the selection position used by the author
$ws.Range("E21").Select()
